$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "220.00026959066784" -as [double]
$ws.Range("B2").Value = "250.00072582102888" -as [double]
$ws.Range("C2").Value = "10.864927139321171" -as [double]
$ws.Range("D2").Value = "0.9999997243416573" -as [double]
$ws.Range("E2").Value = "-0.0010412608900588801" -as [double]
$ws.Range("F2").Value = "0.0003867540448786701" -as [double]
$ws.Range("G2").Value = "-1.0487571689399101e-16" -as [double]

$ws.Range("A3").Value = "219.99273067256067" -as [double]
$ws.Range("B3").Value = "250.00143886042096" -as [double]
$ws.Range("C3").Value = "11.65047097814216" -as [double]
$ws.Range("D3").Value = "0.9999439909785484" -as [double]
$ws.Range("E3").Value = "-0.0010652922919146992" -as [double]
$ws.Range("F3").Value = "-0.01123764832025641" -as [double]
$ws.Range("G3").Value = "-2.6498588398049186e-6" -as [double]

$ws.Range("A4").Value = "219.98482525636933" -as [double]
$ws.Range("B4").Value = "250.0021779609179" -as [double]
$ws.Range("C4").Value = "12.363793598199235" -as [double]
$ws.Range("D4").Value = "0.9999332887766554" -as [double]
$ws.Range("E4").Value = "-0.001145912195236593" -as [double]
$ws.Range("F4").Value = "-0.012221470947576212" -as [double]
$ws.Range("G4").Value = "-3.0199053615329836e-6" -as [double]

$ws.Range("A5").Value = "219.97664596688688" -as [double]
$ws.Range("B5").Value = "250.002942479048" -as [double]
$ws.Range("C5").Value = "13.011525438739707" -as [double]
$ws.Range("D5").Value = "0.9999231552673236" -as [double]
$ws.Range("E5").Value = "-0.0012273724887350728" -as [double]
$ws.Range("F5").Value = "-0.01308612498073136" -as [double]
$ws.Range("G5").Value = "-2.8583672499424103e-6" -as [double]

$ws.Range("A6").Value = "219.96449995234218" -as [double]
$ws.Range("B6").Value = "249.95724406499687" -as [double]
$ws.Range("C6").Value = "13.61950846179564" -as [double]
$ws.Range("D6").Value = "0.9970140514007063" -as [double]
$ws.Range("E6").Value = "0.07463022813286826" -as [double]
$ws.Range("F6").Value = "-0.019834565490675724" -as [double]
$ws.Range("G6").Value = "1.1050578884936881e-5" -as [double]

$ws.Range("A7").Value = "219.92773771811352" -as [double]
$ws.Range("B7").Value = "249.91059641175895" -as [double]
$ws.Range("C7").Value = "14.184271806018279" -as [double]
$ws.Range("D7").Value = "0.9950620863234991" -as [double]
$ws.Range("E7").Value = "0.0779421975263948" -as [double]
$ws.Range("F7").Value = "-0.06142401491881902" -as [double]
$ws.Range("G7").Value = "2.4695313743101066e-5" -as [double]

$ws.Range("A8").Value = "219.88961695483405" -as [double]
$ws.Range("B8").Value = "249.86222542773598" -as [double]
$ws.Range("C8").Value = "14.701036800776686" -as [double]
$ws.Range("D8").Value = "0.9943953189752971" -as [double]
$ws.Range("E8").Value = "0.08301959626795787" -as [double]
$ws.Range("F8").Value = "-0.06542597445942969" -as [double]
$ws.Range("G8").Value = "2.7287183986506357e-5" -as [double]

$ws.Range("A9").Value = "219.85011030953697" -as [double]
$ws.Range("B9").Value = "249.81209652400972" -as [double]
$ws.Range("C9").Value = "15.174524982587489" -as [double]
$ws.Range("D9").Value = "0.993671690646245" -as [double]
$ws.Range("E9").Value = "0.08819537115938278" -as [double]
$ws.Range("F9").Value = "-0.06950562808979209" -as [double]
$ws.Range("G9").Value = "3.05186162253215e-5" -as [double]

$ws.Range("A10").Value = "219.8091936407595" -as [double]
$ws.Range("B10").Value = "249.76017934668008" -as [double]
$ws.Range("C10").Value = "15.609023440746506" -as [double]
$ws.Range("D10").Value = "0.9928915770009736" -as [double]
$ws.Range("E10").Value = "0.09344958080035048" -as [double]
$ws.Range("F10").Value = "-0.07364748446434037" -as [double]
$ws.Range("G10").Value = "3.535272598553793e-5" -as [double]

$ws.Range("A11").Value = "219.73954619148384" -as [double]
$ws.Range("B11").Value = "249.7093907453672" -as [double]
$ws.Range("C11").Value = "16.036125602006088" -as [double]
$ws.Range("D11").Value = "0.9878272011837405" -as [double]
$ws.Range("E11").Value = "0.09159283523281816" -as [double]
$ws.Range("F11").Value = "-0.12560094001206648" -as [double]
$ws.Range("G11").Value = "0.00023278098363854814" -as [double]

$ws.Range("A12").Value = "219.6679373819163" -as [double]
$ws.Range("B12").Value = "249.65727047873946" -as [double]
$ws.Range("C12").Value = "16.431731619286747" -as [double]
$ws.Range("D12").Value = "0.9866532613061523" -as [double]
$ws.Range("E12").Value = "0.09575575674976064" -as [double]
$ws.Range("F12").Value = "-0.13155928247910997" -as [double]
$ws.Range("G12").Value = "0.000758463146040114" -as [double]

$ws.Range("A13").Value = "219.5936353126526" -as [double]
$ws.Range("B13").Value = "249.6031638789197" -as [double]
$ws.Range("C13").Value = "16.79915182549215" -as [double]
$ws.Range("D13").Value = "0.9851200694359724" -as [double]
$ws.Range("E13").Value = "0.10109282957644561" -as [double]
$ws.Range("F13").Value = "-0.13882107301457708" -as [double]
$ws.Range("G13").Value = "0.0006753678387950495" -as [double]

$ws.Range("A14").Value = "219.51704466106307" -as [double]
$ws.Range("B14").Value = "249.54739052195106" -as [double]
$ws.Range("C14").Value = "17.141267464514595" -as [double]
$ws.Range("D14").Value = "0.9836858700548914" -as [double]
$ws.Range("E14").Value = "0.10581198030765823" -as [double]
$ws.Range("F14").Value = "-0.14529834759271987" -as [double]
$ws.Range("G14").Value = "0.0007052212097034686" -as [double]

$ws.Range("A15").Value = "219.43825361630363" -as [double]
$ws.Range("B15").Value = "249.4898368554855" -as [double]
$ws.Range("C15").Value = "17.460648132730213" -as [double]
$ws.Range("D15").Value = "0.9821956661599742" -as [double]
$ws.Range("E15").Value = "0.1107167547391844" -as [double]
$ws.Range("F15").Value = "-0.15155674076719589" -as [double]
$ws.Range("G15").Value = "-0.00010830276304162355" -as [double]

$ws.Range("A16").Value = "219.35629592404143" -as [double]
$ws.Range("B16").Value = "249.4298201276793" -as [double]
$ws.Range("C16").Value = "17.76103615358067" -as [double]
$ws.Range("D16").Value = "0.9802270793194887" -as [double]
$ws.Range("E16").Value = "0.11678155688779235" -as [double]
$ws.Range("F16").Value = "-0.15945745455814336" -as [double]
$ws.Range("G16").Value = "-0.0008415078396036129" -as [double]

$ws.Range("A17").Value = "219.27155215296693" -as [double]
$ws.Range("B17").Value = "249.36791777525917" -as [double]
$ws.Range("C17").Value = "18.044254614897454" -as [double]
$ws.Range("D17").Value = "0.978428224325585" -as [double]
$ws.Range("E17").Value = "0.12171379750330237" -as [double]
$ws.Range("F17").Value = "-0.1665963224076738" -as [double]
$ws.Range("G17").Value = "-0.0001500672239981155" -as [double]

$ws.Range("A18").Value = "219.18373270558683" -as [double]
$ws.Range("B18").Value = "249.30378772476186" -as [double]
$ws.Range("C18").Value = "18.313073776788997" -as [double]
$ws.Range("D18").Value = "0.9764019892059159" -as [double]
$ws.Range("E18").Value = "0.12718615280986426" -as [double]
$ws.Range("F18").Value = "-0.17412788308371938" -as [double]
$ws.Range("G18").Value = "-0.000126314245383111" -as [double]

$ws.Range("A19").Value = "219.09240583672084" -as [double]
$ws.Range("B19").Value = "249.237084917124" -as [double]
$ws.Range("C19").Value = "18.570781151422633" -as [double]
$ws.Range("D19").Value = "0.9740983017081197" -as [double]
$ws.Range("E19").Value = "0.13308863671314408" -as [double]
$ws.Range("F19").Value = "-0.18224852323641422" -as [double]
$ws.Range("G19").Value = "-0.00020590113756349736" -as [double]

$ws.Range("A20").Value = "218.99805314206876" -as [double]
$ws.Range("B20").Value = "249.16816036031378" -as [double]
$ws.Range("C20").Value = "18.818419429590847" -as [double]
$ws.Range("D20").Value = "0.9719910328905776" -as [double]
$ws.Range("E20").Value = "0.13830638275267207" -as [double]
$ws.Range("F20").Value = "-0.18937242504301446" -as [double]
$ws.Range("G20").Value = "-0.0002832467056612076" -as [double]

$ws.Range("A21").Value = "218.9004431336227" -as [double]
$ws.Range("B21").Value = "249.09684267918848" -as [double]
$ws.Range("C21").Value = "19.058398769305143" -as [double]
$ws.Range("D21").Value = "0.9697294460284326" -as [double]
$ws.Range("E21").Value = "0.14368148135113584" -as [double]
$ws.Range("F21").Value = "-0.19670736604813496" -as [double]
$ws.Range("G21").Value = "-0.0003824864448885051" -as [double]

$ws.Range("A22").Value = "218.7991381217443" -as [double]
$ws.Range("B22").Value = "249.02282470289032" -as [double]
$ws.Range("C22").Value = "19.29367728205434" -as [double]
$ws.Range("D22").Value = "0.9672127206048771" -as [double]
$ws.Range("E22").Value = "0.1494309480434845" -as [double]
$ws.Range("F22").Value = "-0.20454993919117218" -as [double]
$ws.Range("G22").Value = "-0.0004253820547333859" -as [double]

$ws.Range("A23").Value = "218.6940412213525" -as [double]
$ws.Range("B23").Value = "248.94602157809973" -as [double]
$ws.Range("C23").Value = "19.526499669690068" -as [double]
$ws.Range("D23").Value = "0.9646205428510412" -as [double]
$ws.Range("E23").Value = "0.1551065403605994" -as [double]
$ws.Range("F23").Value = "-0.21228604447546987" -as [double]
$ws.Range("G23").Value = "-0.0005439700505425281" -as [double]

$ws.Range("A24").Value = "218.5834715295162" -as [double]
$ws.Range("B24").Value = "248.86343590137713" -as [double]
$ws.Range("C24").Value = "19.764696659492888" -as [double]
$ws.Range("D24").Value = "0.9605501480147939" -as [double]
$ws.Range("E24").Value = "0.16585870930781635" -as [double]
$ws.Range("F24").Value = "-0.2220583877815519" -as [double]
$ws.Range("G24").Value = "-0.00437740028453201" -as [double]

$ws.Range("A25").Value = "218.47715381161638" -as [double]
$ws.Range("B25").Value = "248.78292683252894" -as [double]
$ws.Range("C25").Value = "19.938735321279754" -as [double]
$ws.Range("D25").Value = "0.9597028230279327" -as [double]
$ws.Range("E25").Value = "0.1690233745090312" -as [double]
$ws.Range("F25").Value = "-0.2232051291347308" -as [double]
$ws.Range("G25").Value = "-0.007367470239883023" -as [double]

$ws.Range("A26").Value = "218.37177475810373" -as [double]
$ws.Range("B26").Value = "248.75622717838903" -as [double]
$ws.Range("C26").Value = "19.99999999936089" -as [double]
$ws.Range("D26").Value = "0.9670473932501742" -as [double]
$ws.Range("E26").Value = "0.0611917894232326" -as [double]
$ws.Range("F26").Value = "-0.24142079628490096" -as [double]
$ws.Range("G26").Value = "-0.04918256278121028" -as [double]

$ws.Range("A27").Value = "218.27369679799597" -as [double]
$ws.Range("B27").Value = "248.7370425735193" -as [double]
$ws.Range("C27").Value = "19.99999999948656" -as [double]
$ws.Range("D27").Value = "0.9700268010717867" -as [double]
$ws.Range("E27").Value = "0.04624375911219666" -as [double]
$ws.Range("F27").Value = "-0.23633359705444645" -as [double]
$ws.Range("G27").Value = "-0.02751738408388585" -as [double]

$ws.Range("A28").Value = "218.17561967543168" -as [double]
$ws.Range("B28").Value = "248.71801138394537" -as [double]
$ws.Range("C28").Value = "19.999999999480522" -as [double]
$ws.Range("D28").Value = "0.970047135924364" -as [double]
$ws.Range("E28").Value = "0.045874550605186235" -as [double]
$ws.Range("F28").Value = "-0.23633421033063828" -as [double]
$ws.Range("G28").Value = "-0.027420979976970277" -as [double]

$ws.Range("A29").Value = "218.0775434490379" -as [double]
$ws.Range("B29").Value = "248.69913643830455" -as [double]
$ws.Range("C29").Value = "19.999999999456506" -as [double]
$ws.Range("D29").Value = "0.9700677931555529" -as [double]
$ws.Range("E29").Value = "0.04549852424099319" -as [double]
$ws.Range("F29").Value = "-0.23633472903153888" -as [double]
$ws.Range("G29").Value = "-0.027320095940505197" -as [double]

$ws.Range("A30").Value = "217.97946817994898" -as [double]
$ws.Range("B30").Value = "248.6804206497972" -as [double]
$ws.Range("C30").Value = "19.999999999479446" -as [double]
$ws.Range("D30").Value = "0.9700887795180559" -as [double]
$ws.Range("E30").Value = "0.04511547591503826" -as [double]
$ws.Range("F30").Value = "-0.2363351483561987" -as [double]
$ws.Range("G30").Value = "-0.027214541252565897" -as [double]

$ws.Range("A31").Value = "217.8813939317542" -as [double]
$ws.Range("B31").Value = "248.66186701948467" -as [double]
$ws.Range("C31").Value = "19.999999999485386" -as [double]
$ws.Range("D31").Value = "0.9701101018036861" -as [double]
$ws.Range("E31").Value = "0.04472519349626339" -as [double]
$ws.Range("F31").Value = "-0.2363354633750834" -as [double]
$ws.Range("G31").Value = "-0.027104115938104342" -as [double]

$ws.Range("A32").Value = "217.78332077065375" -as [double]
$ws.Range("B32").Value = "248.643478639795" -as [double]
$ws.Range("C32").Value = "19.99999999948542" -as [double]
$ws.Range("D32").Value = "0.9701317668554412" -as [double]
$ws.Range("E32").Value = "0.04432745640182336" -as [double]
$ws.Range("F32").Value = "-0.23633566893326208" -as [double]
$ws.Range("G32").Value = "-0.026988610219126102" -as [double]

$ws.Range("A33").Value = "217.6852487655373" -as [double]
$ws.Range("B33").Value = "248.62525869817247" -as [double]
$ws.Range("C33").Value = "19.999999999547505" -as [double]
$ws.Range("D33").Value = "0.9701537815408353" -as [double]
$ws.Range("E33").Value = "0.04392203521981299" -as [double]
$ws.Range("F33").Value = "-0.23633575968294565" -as [double]
$ws.Range("G33").Value = "-0.026867803988810612" -as [double]

$ws.Range("A34").Value = "217.58717798805597" -as [double]
$ws.Range("B34").Value = "248.60721048089178" -as [double]
$ws.Range("C34").Value = "19.999999999576055" -as [double]
$ws.Range("D34").Value = "0.9701761527308134" -as [double]
$ws.Range("E34").Value = "0.04350869130727064" -as [double]
$ws.Range("F34").Value = "-0.2363357300837116" -as [double]
$ws.Range("G34").Value = "-0.026741466261453677" -as [double]

$ws.Range("A35").Value = "217.48910851268698" -as [double]
$ws.Range("B35").Value = "248.5893373770382" -as [double]
$ws.Range("C35").Value = "19.999999999546677" -as [double]
$ws.Range("D35").Value = "0.9701988872744668" -as [double]
$ws.Range("E35").Value = "0.043087176385877425" -as [double]
$ws.Range("F35").Value = "-0.23633557440215588" -as [double]
$ws.Range("G35").Value = "-0.02660935461909175" -as [double]

